# Mark petroleum and heavy or residual fuel oil plants as peakers
$wb = $excel.ActiveWorkbook

# "BPaFF-BITPTaP" = Boolean Is This Plant Type a Peaker
$wsPeaker = $wb.Worksheets.Item("BPaFF-BITPTaP")

# Flag petroleum (row 11) as a peaker plant type.
$wsPeaker.Range("B11").Value = 1

# "heavy or residual fuel oil" (row 16) previously mirrored "petroleum" (B11) via
# a formula; it is now set directly as a peaker flag (literal 1), no longer linked
# to B11 by formula.
$wsPeaker.Range("B16").Formula = "1"

# Move the active selection to reflect where the edit was made, then
# restore "About" as the active/displayed sheet (unchanged from before).
$wsPeaker.Range("B12").Select()
$wb.Worksheets.Item("About").Activate()

$wb.Save()
